$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Formula = "=A11+B11"
$ws.Range("C11").NumberFormat = $ws.Range("C10").NumberFormat

$ws.Range("C2:C11").Select() | Out-Null
